$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2)
}

# 1. "Ano 2080," -> "Ano 2141,"
Replace-Text "Ano 2080," "Ano 2141,"

# 2. "pois o mesmo deveria ter que enfrentar" -> "pois o mesmo teria de enfrentar"
Replace-Text "pois o mesmo deveria ter que enfrentar" "pois o mesmo teria de enfrentar"

# 3. "da raça XXX bem" -> "da raça Zetagrin bem"
Replace-Text "da raça XXX bem" "da raça Zetagrin bem"

# 4. "Ele nasceu em 2069 e" -> "Ele nasceu em 2126 e"
Replace-Text "Ele nasceu em 2069 e" "Ele nasceu em 2126 e"

# 5. "despesas de casa, a partir de então ele começou" -> "despesas de casa, a partir de então, ele começou"
Replace-Text "despesas de casa, a partir de então ele começou" "despesas de casa, a partir de então, ele começou"

# 6. "do universo, afinal somente ele tinha" -> "do universo, afinal, somente ele tinha"
Replace-Text "do universo, afinal somente ele tinha" "do universo, afinal, somente ele tinha"

# 7. "criada por ele mesmo para poder respirar em outros planetas e suas armas"
#    -> "criada por ele mesmo, para poder respirar em outros planetas, e suas armas"
Replace-Text "criada por ele mesmo para poder respirar em outros planetas e suas armas" "criada por ele mesmo, para poder respirar em outros planetas, e suas armas"

# 8. "Fase 1 – Fábrica de Vulcões" -> "Fase 1 – Indústria de Vulcões"
Replace-Text "Fase 1 – Fábrica de Vulcões" "Fase 1 – Indústria de Vulcões"

# 9. "iniciadas pela Moldavia em 2078" -> "iniciadas pela Moondavia em 2078"
Replace-Text "iniciadas pela Moldavia em 2078" "iniciadas pela Moondavia em 2078"

# 10. "pesquisadores Zhargox e Yorkaxig, seus bisnetos" -> "pesquisadores Zhargox e Yorkaxig. Seus bisnetos"
Replace-Text "pesquisadores Zhargox e Yorkaxig, seus bisnetos" "pesquisadores Zhargox e Yorkaxig. Seus bisnetos"

# 11. "é um dos planetas recentemente descoberto com apenas 7 anos"
#     -> "é um planeta recentemente descoberto com apenas sete anos"
Replace-Text "é um dos planetas recentemente descoberto com apenas 7 anos" "é um planeta recentemente descoberto com apenas sete anos"

# 12. "Algumas pessoas já fugiram, dizendo que este planeta está corrompido por alguma forma maléfica."
#     -> "Algumas pessoas já fugiram dizendo que este planeta está corrompido por alguma forma maléfica e outras ainda procuram por alguma solução para o que está acontecendo, mas não sabem até quando tudo permanecerá desta maneira."
Replace-Text "Algumas pessoas já fugiram, dizendo que este planeta está corrompido por alguma forma maléfica." "Algumas pessoas já fugiram dizendo que este planeta está corrompido por alguma forma maléfica e outras ainda procuram por alguma solução para o que está acontecendo, mas não sabem até quando tudo permanecerá desta maneira."

# 13. Move the _GoBack bookmark from after "Atualmente ... certo." to right
#     after "...desta maneira" (before the final period) in the last paragraph.
try {
    $d.Bookmarks.Item("_GoBack").Delete()
} catch {
}

$bmRange = $d.Content
$bmRange.Find.Execute("desta maneira", $true, $false, $false, $false, $false,
                       $true, 1, $false, "", 0) | Out-Null
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange)
